$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: change health from GREEN to RED, rename project to "o3"
$ws.Range("C23").Value = "o3"
$ws.Range("G23").Value = "RED"
$ws.Range("O23").Value = 0.0

# Row 24: rename project to "TEST", clear health (blank = counted as green),
# change service category to "Threat Management"
$ws.Range("C24").Value = "TEST"
$ws.Range("G24").ClearContents()
$ws.Range("I24").Value = "Threat Management"
$ws.Range("O24").Value = 0.5

# Remove the now-unused extra test rows (25-35) entirely
$ws.Range("B25:P35").Clear()

# Shrink the table to match the reduced data range
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:O24"))
